# Swap the deck's theme palette: the slide master currently uses the
# "Integral" color scheme (ppt/theme/theme1.xml); the commit replaces it
# with the stock "Office Theme" color scheme (the palette that used to
# live on the notes master's theme, ppt/theme/theme2.xml).
#
# PowerPoint exposes the twelve theme colors (dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink -- in that fixed index order) through
# Master.Theme.ThemeColorScheme.Colors(n).RGB, so we drive the swap the
# same way a user would from Design > Variants > Colors > Customize
# Colors, setting each of the twelve slots on the slide master's theme.

$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$tcs = $m.Theme.ThemeColorScheme

function Set-ThemeColor($index, $hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    $tcs.Colors($index).RGB = $r + ($g * 256) + ($b * 65536)
}

# Office Theme palette (was previously theme2.xml / notes master theme).
Set-ThemeColor 1  "000000"   # dk1
Set-ThemeColor 2  "FFFFFF"   # lt1
Set-ThemeColor 3  "44546A"   # dk2
Set-ThemeColor 4  "E7E6E6"   # lt2
Set-ThemeColor 5  "5B9BD5"   # accent1
Set-ThemeColor 6  "ED7D31"   # accent2
Set-ThemeColor 7  "A5A5A5"   # accent3
Set-ThemeColor 8  "FFC000"   # accent4
Set-ThemeColor 9  "4472C4"   # accent5
Set-ThemeColor 10 "70AD47"   # accent6
Set-ThemeColor 11 "0563C1"   # hlink
Set-ThemeColor 12 "954F72"   # folHlink
